$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '27.935.64'
$ws.Range('E2').Value = '  +1.42%  '

$ws.Range('D3').Value = '1.642.92'
$ws.Range('E3').Value = '  +1.19%  '

$ws.Range('E4').Value = '  -0.01%  '

Set-TextValue $ws.Range('D5') '213.47'

$ws.Range('E6').Value = '  -0.14%  '

$ws.Range('E7').Value = '  -0.03%  '

Set-TextValue $ws.Range('D8') '23.86'
$ws.Range('E8').Value = '  +2.66%  '

$ws.Range('E9').Value = '  +0.52%  '

$ws.Range('E10').Value = '  +0.88%  '

$ws.Range('E11').Value = '  -0.92%  '

$ws.Range('D12').Value = '1.876.18'
$ws.Range('E12').Value = '  +1.21%  '

$ws.Range('D13').Value = '1.648.90'
$ws.Range('E13').Value = '  +1.56%  '

Set-TextValue $ws.Range('D14') '0.575'
$ws.Range('E14').Value = '  +4.64%  '

$ws.Range('E15').Value = '  +0.92%  '

$ws.Range('E16').Value = '  +1.13%  '

$ws.Range('D17').Value = '27.916.49'
$ws.Range('E17').Value = '  +1.44%  '

Set-TextValue $ws.Range('D18') '230.63'

$ws.Range('E19').Value = '  +1.00%  '

Set-TextValue $ws.Range('D20') '7.63'
$ws.Range('E20').Value = '  +1.41%  '

$ws.Range('E21').Value = '  -0.04%  '

Set-TextValue $ws.Range('D22') '11.12'
$ws.Range('E22').Value = '  +7.45%  '

Set-TextValue $ws.Range('D23') '4.40'
$ws.Range('E23').Value = '  +1.63%  '

Set-TextValue $ws.Range('D24') '2.06'
$ws.Range('E24').Value = '  -0.69%  '

Set-TextValue $ws.Range('D25') '152.92'
$ws.Range('E25').Value = '  +2.78%  '

$ws.Range('E26').Value = '  +0.86%  '

$ws.Range('E27').Value = '  +0.81%  '

Set-TextValue $ws.Range('D28') '15.72'
$ws.Range('E28').Value = '  +1.12%  '

$ws.Range('E30').Value = '  +1.09%  '

$ws.Range('E31').Value = '  +0.53%  '

Set-TextValue $ws.Range('D32') '3.33'
$ws.Range('E32').Value = '  +2.03%  '

$ws.Range('D33').Value = '1.424.90'
$ws.Range('E33').Value = '  -2.94%  '

Set-TextValue $ws.Range('D34') '3.11'
$ws.Range('E34').Value = '  +2.01%  '

$ws.Range('E35').Value = '  +2.13%  '

Set-TextValue $ws.Range('D36') '2.34'

Set-TextValue $ws.Range('D37') '0.890'
$ws.Range('E37').Value = '  +1.95%  '

$ws.Range('E38').Value = '  -0.45%  '

$ws.Range('E39').Value = '  +0.97%  '

Set-TextValue $ws.Range('D40') '0.558'
$ws.Range('E40').Value = '  +0.55%  '

$ws.Range('E41').Value = '  +2.48%  '

$ws.Range('E42').Value = '  -0.05%  '

Set-TextValue $ws.Range('D43') '67.23'
$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('E44').Value = '  +0.55%  '

$ws.Range('E45').Value = '  +3.27%  '

Set-TextValue $ws.Range('D46') '1.81'
$ws.Range('E46').Value = '  +2.98%  '

$ws.Range('E47').Value = '  -0.03%  '

$ws.Range('D48').Value = '1.784.63'
$ws.Range('E48').Value = '  +1.14%  '

Set-TextValue $ws.Range('D49') '88.99'
$ws.Range('E49').Value = '  +2.03%  '

$ws.Range('E50').Value = '  +1.15%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D51') '0.0506'
$ws.Range('E51').Value = '  +0.50%  '
